$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Hello and Welcome, `nMy Name is SAM your Search Assistant Manager.`nHow may I help you?"
